# Add the missing weekly price record for "Vega Modelo de Temuco" - Piña.
# A new data row is inserted above the existing row 447, shifting every
# subsequent row down by one (old row 447 becomes 448, ..., old row 473
# becomes 474), and the new row is populated with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 447 (pushes 447..473 down to 448..474).
$ws.Rows.Item(447).Insert()

$ws.Cells.Item(447,1).Value  = 10
$ws.Cells.Item(447,2).Value  = 'Vega Modelo de Temuco'
$ws.Cells.Item(447,3).Value  = 'La Araucanía'
$ws.Cells.Item(447,4).Value  = 44706
$ws.Cells.Item(447,5).Value  = 9
$ws.Cells.Item(447,6).Value  = 'Fruta'
$ws.Cells.Item(447,7).Value  = 100108
$ws.Cells.Item(447,8).Value  = 'Tropicales y subtropicales'
$ws.Cells.Item(447,9).Value  = 100108005
$ws.Cells.Item(447,10).Value = 'Piña'
$ws.Cells.Item(447,11).Value = 'Caramelo'
$ws.Cells.Item(447,12).Value = 'Primera'
$ws.Cells.Item(447,13).Value = 75
$ws.Cells.Item(447,14).Value = 20000
$ws.Cells.Item(447,15).Value = 20000
$ws.Cells.Item(447,16).Value = 20000
$ws.Cells.Item(447,17).Value = '$/caja 12 unidades'
$ws.Cells.Item(447,18).Value = 'Ecuador'
$ws.Cells.Item(447,19).Value = 1667
$ws.Cells.Item(447,20).Value = 12
